# match_environment.xlsx now ships with a short, environment-relative DBs
# path instead of the hard-coded path from Pavel's machine, so CSmatchSht
# can fall back to this default when its new optional parameter is omitted.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "C:\DBs\"

# Leave the cursor parked on the cell that was just edited.
$ws.Range("B1").Select() | Out-Null
